$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.721.88'
$ws.Range("E2").Value = '  -0.19%  '

# Row 3
$ws.Range("D3").Value = '2.540.20'
$ws.Range("E3").Value = '  -0.29%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.63'
$ws.Range("E5").Value = '  +0.46%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.52'
$ws.Range("E6").Value = '  +1.98%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.566'
$ws.Range("E7").Value = '  -0.95%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("E9").Value = '  -1.60%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.47'
$ws.Range("E10").Value = '  -0.77%  '

# Row 11
$ws.Range("E11").Value = '  -0.31%  '

# Row 12
$ws.Range("E12").Value = '  -1.37%  '

# Row 13
$ws.Range("E13").Value = '  +1.02%  '

# Row 14
$ws.Range("D14").Value = '2.936.88'
$ws.Range("E14").Value = '  -0.20%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.578.07'
$ws.Range("E15").Value = '  +0.04%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.38'
$ws.Range("E16").Value = '  -3.58%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.817'
$ws.Range("E17").Value = '  -2.74%  '

# Row 18
$ws.Range("D18").Value = '42.730.55'
$ws.Range("E18").Value = '  -0.24%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.73'
$ws.Range("E19").Value = '  -0.26%  '

# Row 20
$ws.Range("E20").Value = '  -0.57%  '

# Row 21
$ws.Range("E21").Value = '  -0.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.98'
$ws.Range("E22").Value = '  +0.94%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.53'

# Row 24
$ws.Range("E24").Value = '  -1.57%  '

# Row 25
$ws.Range("E25").Value = '  -1.71%  '

# Row 26
$ws.Range("E26").Value = '  +0.12%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.64'
$ws.Range("E27").Value = '  -3.81%  '

# Row 28
$ws.Range("E28").Value = '  -1.05%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.16'
$ws.Range("E29").Value = '  -0.07%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.45'
$ws.Range("E30").Value = '  -4.59%  '

# Row 31
$ws.Range("E31").Value = '  +2.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '157.91'
$ws.Range("E32").Value = '  -0.09%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  +5.73%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.67'
$ws.Range("E34").Value = '  +1.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.16'
$ws.Range("E36").Value = '  -4.16%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.02'
$ws.Range("E37").Value = '  -1.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.98'
$ws.Range("E38").Value = '  -5.36%  '

# Row 39
$ws.Range("E39").Value = '  -0.38%  '

# Row 40
$ws.Range("E40").Value = '  -0.30%  '

# Row 41
$ws.Range("E41").Value = '  +1.12%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.90'

# Row 43
$ws.Range("E43").Value = '  +0.17%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0299'
$ws.Range("E44").Value = '  -0.39%  '

# Row 45
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.28'
$ws.Range("E45").Value = '  +1.88%  '

# Row 46
$ws.Range("D46").Value = '1.999.38'
$ws.Range("E46").Value = '  +0.11%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.11'
$ws.Range("E47").Value = '  +0.15%  '

# Row 48
$ws.Range("D48").Value = '2.782.08'
$ws.Range("E48").Value = '  -0.18%  '

# Row 49
$ws.Range("E49").Value = '  +0.10%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '79.92'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.44'
$ws.Range("E51").Value = '  -1.57%  '
